$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 (sldId 256): "CustomShape 2" - update names/date line
#   "Irina & Rao"  -> "Irina & " + "Bea"
#   "Hilary 2024"  -> "March" + " 2024"
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange

$p1para1 = $tr1.Paragraphs(1, 1)
$rao = $p1para1.Characters(9, 3)
$rao.Text = "Bea"

$p1para2 = $tr1.Paragraphs(2, 1)
$hilary = $p1para2.Characters(1, 6)
$hilary.Text = ""
$p1para2b = $tr1.Paragraphs(2, 1)
$p1para2b.InsertBefore("March") | Out-Null

# ---------------------------------------------------------------------------
# Slide 2 (presentation order, "CustomShape 2" with id 229 -> 2):
#   Cut + Paste back in place so PowerPoint re-numbers the shape id (229 -> 2)
#   and regenerates its creation id, matching the recorded change.
#   Then apply the same name/date edits as slide 1, plus the Michaelmas/2023
#   -> March/2024 edits.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2old = $s2.Shapes.Item(2)
$sh2old.Cut()
$sh2 = $s2.Shapes.Paste().Item(1)

$tr2 = $sh2.TextFrame.TextRange

$p2para1 = $tr2.Paragraphs(1, 1)
$rao2 = $p2para1.Characters(9, 3)
$rao2.Text = "Bea"

$p2para2 = $tr2.Paragraphs(2, 1)
$michaelmas = $p2para2.Characters(1, 10)
$michaelmas.Text = ""
$p2para2b = $tr2.Paragraphs(2, 1)
$p2para2b.InsertBefore("March") | Out-Null

$p2para2c = $tr2.Paragraphs(2, 1)
$year2 = $p2para2c.Characters(6, 5)
$year2.Text = " 2024"

# ---------------------------------------------------------------------------
# Slide 3 (presentation order): collapse the GitHub repository link runs
#   "https://github.com/sraorao" + "/MSD_R_course_HT2024" + " "
#   -> single run "https://github.com/Chelysheva/MSD_R_course_lateMarch2024"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(1)
$tr3 = $sh3.TextFrame.TextRange

$p3para2 = $tr3.Paragraphs(2, 1)
$link = $p3para2.Characters(19, 26)
$link.Text = "https://github.com/Chelysheva/MSD_R_course_lateMarch2024"

$p3para2b = $tr3.Paragraphs(2, 1)
$tail = $p3para2b.Characters(75, 21)
$tail.Text = ""
